$wb = $excel.ActiveWorkbook

# Rename the "SwateTemplateMetadata" sheet to "isa_template"
$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

# Update the selection on that sheet from C13 to B14
$metaSheet.Activate()
$metaSheet.Range("B14").Select()
